# The presentation's design/theme (backing ppt/theme/theme1.xml, used by the
# single slide master) is changed from the "Integral" (Red Violet) theme to
# the built-in "Office Theme" colour scheme:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72
#
# PowerPoint's ColorScheme.Colors(i) indexer walks the theme's 12 colour
# slots in the canonical order: dk1, lt1, dk2, lt2, accent1..accent6,
# hlink, folHlink. RGB values are the usual COM BGR-packed long.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

$officeThemeRgb = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $officeThemeRgb.Count; $i++) {
    $cs.Colors($i).RGB = $officeThemeRgb[$i - 1]
}
